$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sr.No. 1) - updated values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "C11CF43403DA"
$ws.Range("C2").Value = "Epson L382 Printer"
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = "TZSTZ01"
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = "NA"
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 5

# Row 3 (Sr.No. 2) - new/updated values
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "C11CF43403DA"
$ws.Range("C3").Value = "Epson L382 Printer"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = "TZSTZ01"
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 10
$ws.Range("L3").Value = "NA"
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = 5

# Clear the old TaxCode-specific formatting from G2/G3 (now plain numeric cells)
$ws.Range("G2:G3").ClearFormats()

# Apply the "code/description" look (small dark-blue Arial, centered, medium white box border)
# to the ItemCode and Description cells for both rows
$styleRng1 = $ws.Range("B2:C3")
$styleRng1.Font.Name = "Arial"
$styleRng1.Font.Size = 7
$styleRng1.Font.Color = 5975850
$styleRng1.HorizontalAlignment = -4108
$styleRng1.VerticalAlignment = -4160
$styleRng1.Borders.LineStyle = 1
$styleRng1.Borders.Weight = -4138
$styleRng1.Borders.Color = 16777215

# ... and to the Warehouse cells for both rows
$styleRng2 = $ws.Range("I2:I3")
$styleRng2.Font.Name = "Arial"
$styleRng2.Font.Size = 7
$styleRng2.Font.Color = 5975850
$styleRng2.HorizontalAlignment = -4108
$styleRng2.VerticalAlignment = -4160
$styleRng2.Borders.LineStyle = 1
$styleRng2.Borders.Weight = -4138
$styleRng2.Borders.Color = 16777215
